# Auto commit at 2025-09-04  9:07:35.07
#
# Updates the "Metrics" sheet's daily figures (B2:B13). The "today" sheet
# (codeName Sheet4) pulls these via =Metrics!Bn formulas (and E/F columns
# derive from those), so it recalculates automatically. Also refreshes the
# active-cell/active-sheet selection state and the "today" sheet's page
# setup.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metrics sheet: new daily values
# ---------------------------------------------------------------------
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 45656.729999999996
$metrics.Range("B3").Value  = 37099.360000000001
$metrics.Range("B4").Value  = 14429.03
$metrics.Range("B5").Value  = 1807
$metrics.Range("B6").Value  = 3964907.6099999994
$metrics.Range("B7").Value  = 3364626.8399999994
$metrics.Range("B8").Value  = 1143794.71
$metrics.Range("B9").Value  = 152967
$metrics.Range("B10").Value = 32430231.410999831
$metrics.Range("B11").Value = 19394496.910000004
$metrics.Range("B12").Value = 11425503.600000001
$metrics.Range("B13").Value = 1250594

# Metrics sheet's own selection moved from J12 -> D18, and it is no longer
# the active tab.
$metrics.Range("D18").Select()

# ---------------------------------------------------------------------
# 2. "today" sheet: becomes the active tab, selection moves I14 -> H11,
#    and picks up a page setup (paper size / orientation).
# ---------------------------------------------------------------------
$today = $wb.Worksheets.Item("today")

$today.PageSetup.PaperSize = 9
$today.PageSetup.Orientation = 1

$today.Range("H11").Select()
